# Auto commit at 2026-01-13  7:52:33.90
# Refresh the charging-station report data: rows 2-36 get a new day's
# readings (new terminal names / serial numbers in A & B, new "last
# charge end" timestamps in C, and a new "as-of" timestamp in D), and
# the active selection moves up one row (A12:D36 -> A11:D36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  A="长沙市开福区高岭香江国际城充电站建设项目"; B="303号直流";         C=46034.541192129633; D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=3;  A="长沙市开福区高岭香江国际城充电站建设项目"; B="212号直流";         C=46034.551041666666; D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=4;  A="长沙市开福区高岭香江国际城充电站建设项目"; B="203号直流";         C=46034.554074074076; D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=5;  A="长沙市开福区高岭香江国际城充电站建设项目"; B="107号直流";         C=46034.567893518521; D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=6;  A="长沙市开福区高岭香江国际城充电站建设项目"; B="108号直流";         C=46034.69122685185;  D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=7;  A="长沙市开福区高岭香江国际城充电站建设项目"; B="109号直流";         C=46034.735671296294; D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=8;  A="长沙市开福区高岭香江国际城充电站建设项目"; B="106号直流";         C=46034.740752314814; D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=9;  A="长沙市开福区高岭香江国际城充电站建设项目"; B="309号直流";         C=46034.750787037039; D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=10; A="长沙市开福区高岭香江国际城充电站建设项目"; B="209号直流";         C=46034.771516203706; D=46035.27449074074;  BIsNumericLooking=$false },
    @{ Row=11; A="飞狐四方坪东区充电站";                     B="9176699442100801"; C=46030.706087962964; D=46035.275196759256; BIsNumericLooking=$false },
    @{ Row=12; A="飞狐四方坪西区充电站";                     B="9176699400500102"; C=46033.590543981481; D=46035.275196759256; BIsNumericLooking=$false },
    @{ Row=13; A="飞狐四方坪西区充电站";                     B="9176699400500204"; C=46034.027974537035; D=46035.275196759256; BIsNumericLooking=$false },
    @{ Row=14; A="飞狐四方坪西区充电站";                     B="9176699400501105"; C=46034.100763888891; D=46035.275196759256; BIsNumericLooking=$false },
    @{ Row=15; A="飞狐四方坪西区充电站";                     B="9176699355900101"; C=46034.106539351851; D=46035.275196759256; BIsNumericLooking=$false },
    @{ Row=16; A="飞狐四方坪西区充电站";                     B="9176699400501101"; C=46034.534131944441; D=46035.275196759256; BIsNumericLooking=$false },
    @{ Row=17; A="飞狐四方坪南区充电站";                     B="9176699368200406"; C=46034.539756944447; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=18; A="飞狐四方坪南区充电站";                     B="9176699368200306"; C=46034.539942129632; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=19; A="飞狐四方坪西区充电站";                     B="9176699400500201"; C=46034.540208333332; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=20; A="飞狐四方坪东区充电站";                     B="9176699425700301"; C=46034.546307870369; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=21; A="飞狐四方坪西区充电站";                     B="9176699355900102"; C=46034.550486111111; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=22; A="飞狐四方坪东区充电站";                     B="9176699416300203"; C=46034.553773148145; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=23; A="飞狐四方坪西区充电站";                     B="9176699400500205"; C=46034.557511574072; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=24; A="飞狐四方坪西区充电站";                     B="9176699400500203"; C=46034.560474537036; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=25; A="飞狐四方坪西区充电站";                     B="9176699400501205"; C=46034.576886574076; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=26; A="飞狐四方坪西区充电站";                     B="9176699400500604"; C=46034.583368055559; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=27; A="飞狐四方坪东区充电站";                     B="9176699420300104"; C=46034.585405092592; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=28; A="飞狐四方坪西区充电站";                     B="9176699400501203"; C=46034.595497685186; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=29; A="飞狐四方坪西区充电站";                     B="9176699400501104"; C=46034.601331018515; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=30; A="飞狐四方坪西区充电站";                     B="9176699400501303"; C=46034.614247685182; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=31; A="飞狐四方坪东区充电站";                     B="9176699442100202"; C=46034.618287037039; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=32; A="飞狐四方坪南区充电站";                     B="9176699368200104"; C=46034.635370370372; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=33; A="飞狐四方坪西区充电站";                     B="9176699400500505"; C=46034.663217592592; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=34; A="飞狐四方坪东区充电站";                     B="9176699442101001"; C=46034.687476851854; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=35; A="飞狐四方坪南区充电站";                     B="9176699368200201"; C=46034.764120370368; D=46035.275196759256; BIsNumericLooking=$true },
    @{ Row=36; A="飞狐四方坪南区充电站";                     B="9176699368200203"; C=46034.771064814813; D=46035.275196759256; BIsNumericLooking=$true }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A: always a station-name string (never digit-only), plain assignment is safe.
    $ws.Cells.Item($row, 1).Value = $r.A

    # Column B: terminal name / serial number. For rows whose B cells use the
    # "General" number format, Excel's automatic type-detection would read a
    # long all-digit string as a number and silently round it (IEEE-754 double
    # only carries ~15-16 significant digits, these serials run to 16-17).
    # Force a text entry by flipping the cell to Text format, writing the
    # value, then flipping back to General so the stored style index is left
    # exactly as it was.
    if ($r.BIsNumericLooking) {
        $cell = $ws.Cells.Item($row, 2)
        $cell.NumberFormat = "@"
        $cell.Value = $r.B
        $cell.NumberFormat = "general"
    } else {
        $ws.Cells.Item($row, 2).Value = $r.B
    }

    # Columns C & D: plain date/time serial numbers.
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}

# Selection moves from A12:D36 to A11:D36
$ws.Range("A11:D36").Select()
